$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Remove the duplicated "bold title" paragraph and replace the
#    text of the final (italic) paragraph near the end of the
#    document with the new "Create a feature image..." prompt text.
#    Doing this first (while paragraph numbering still matches the
#    original document) keeps the indices simple.
# ---------------------------------------------------------------

$lastCount = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($lastCount - 1)
Write-Host "Removing duplicate title paragraph: [" $dupTitlePara.Range.Text "]"
$dupTitlePara.Range.Delete()

$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalRange = $finalPara.Range
# Exclude the trailing paragraph mark so only the run text is replaced
# (this preserves the run's existing italic formatting).
$finalTextRange = $d.Range($finalRange.Start, $finalRange.End - 1)
$finalTextRange.Text = "Create a feature image for African Elephant that showcases a happy Maya warrior with glasses on a safari tour, surrounded by the game's iconic animals such as a cheetah, buffalo, and of course, an African elephant. The image should be in a cartoon style that's colorful and vibrant, with elements of the savanna in the background, like trees and grass. Place the Maya warrior front and center, with a big smile on his face and a camera in hand, ready to capture the animals' beauty. The image should also include the game's logo and some of the game's symbols, like the elephant and the playing cards. Make sure the image exudes excitement and adventure to entice players to try out the game."

# ---------------------------------------------------------------
# 2) Insert a new "Meta description" paragraph right after the
#    document's title (Heading1) paragraph.
# ---------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null

$metaPara = $d.Paragraphs.Item(2).Range
$metaPara.Style = "Normal"
$metaPara.Text = "Meta description: Play African Elephant for free and experience an exciting safari-themed slot game. Dive into the vivid savannahs of Africa and win big."

# Make "Meta description" bold, matching the target formatting.
$metaParaRange = $d.Paragraphs.Item(2).Range
$boldRange = $d.Range($metaParaRange.Start, $metaParaRange.Start + 16)
$boldRange.Bold = 1

# Add a leading empty run, mirroring the structural pattern used by
# all of the other body paragraphs in this document.
$collapsedStart = $d.Range($metaParaRange.Start, $metaParaRange.Start)
$collapsedStart.InsertBefore("")

Write-Host "Done."
